# Atualização de bases das ligas, do dia: 20-06-2024 às 20:11
#
# The source data rows for several matches were written to the wrong row
# (their B:AD data - id, teams, scores, odds, etc. - was swapped with an
# adjacent row that shares the same match date). This script restores the
# correct row for each affected match by swapping (or, for one group of
# three rows, cyclically rotating) the data range B:AD between the rows
# involved, while leaving column A (the running row index) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($rowA, $rowB) {
    $rangeA = $ws.Range("B$rowA`:AD$rowA")
    $rangeB = $ws.Range("B$rowB`:AD$rowB")
    $valA = $rangeA.Value2
    $valB = $rangeB.Value2
    $rangeA.Value2 = $valB
    $rangeB.Value2 = $valA
}

# Rotate the data of rows r1 -> r2 -> r3 -> r1, i.e.
#   new(r1) = old(r2), new(r2) = old(r3), new(r3) = old(r1)
function Rotate-Rows($r1, $r2, $r3) {
    $range1 = $ws.Range("B$r1`:AD$r1")
    $range2 = $ws.Range("B$r2`:AD$r2")
    $range3 = $ws.Range("B$r3`:AD$r3")
    $val1 = $range1.Value2
    $val2 = $range2.Value2
    $val3 = $range3.Value2
    $range1.Value2 = $val2
    $range2.Value2 = $val3
    $range3.Value2 = $val1
}

# Simple pairwise swaps
Swap-Rows 18 19
Swap-Rows 105 106
Swap-Rows 117 118
Swap-Rows 133 134
Swap-Rows 150 151
Swap-Rows 164 165
Swap-Rows 197 198
Swap-Rows 203 204
Swap-Rows 210 211
Swap-Rows 218 219

# Three-way rotation
Rotate-Rows 158 159 160
